$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-06-11"

# Update the column header text (shared string) for the "Total" column
$ws.Range("I1").Value = "2022 (through 06-11)"

# Update the updated data values
$ws.Range("I7").Value = 40
$ws.Range("I14").Value = 703
